$d = $word.ActiveDocument
Write-Output ("before paragraphs: " + $d.Paragraphs.Count)
$r = $d.Range(0, $d.Content.End)
$r.InsertXML("<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:r><w:t>ONLY-THIS</w:t></w:r></w:p>")
Write-Output ("after paragraphs: " + $d.Paragraphs.Count)
Write-Output $d.Content.Text
